# Auto-generated Excel COM-interop edit script
# Updates cached market-price / profit figures across all 8 sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per the scheduled-runner refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 570.1667
$ws.Range("J53").Value = 438.33334
$ws.Range("L53").Value = 438.33334
$ws.Range("N53").Value = -1712.33334

$ws.Range("H86").Value = 2449.9167
$ws.Range("I86").Value = 2180
$ws.Range("J86").Value = 2642.7144
$ws.Range("K86").Value = 2180
$ws.Range("L86").Value = 2642.7144
$ws.Range("M86").Value = -1057
$ws.Range("N86").Value = -4888.7144

$ws.Range("H89").Value = 2449.9167
$ws.Range("I89").Value = 2180
$ws.Range("J89").Value = 2642.7144
$ws.Range("K89").Value = 10900
$ws.Range("L89").Value = 13213.572
$ws.Range("M89").Value = -5284
$ws.Range("N89").Value = -24445.572

$ws.Range("H96").Value = 1711.7273
$ws.Range("I96").Value = 1133.2858
$ws.Range("J96").Value = 2724
$ws.Range("K96").Value = 3399.8574
$ws.Range("L96").Value = 8172
$ws.Range("M96").Value = -2026.8574
$ws.Range("N96").Value = -10918

$ws.Range("H99").Value = 823.5
$ws.Range("J99").Value = 1500
$ws.Range("L99").Value = 4500
$ws.Range("N99").Value = -7496

$ws.Range("H129").Value = 2037.1818
$ws.Range("I129").Value = 2037.1818
$ws.Range("K129").Value = 6111.5454
$ws.Range("M129").Value = -1111.5454

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7940151.5
$ws.Range("I32").Value = 7940151.5
$ws.Range("K32").Value = 7940151.5
$ws.Range("M32").Value = -7939864.5

$ws.Range("H45").Value = 2526.3333
$ws.Range("I45").Value = 2567.125
$ws.Range("K45").Value = 2567.125
$ws.Range("M45").Value = -2190.125

$ws.Range("H74").Value = 2689.7742
$ws.Range("I74").Value = 903.65216
$ws.Range("K74").Value = 903.65216
$ws.Range("M74").Value = -29.65215999999998

$ws.Range("H77").Value = 2689.7742
$ws.Range("I77").Value = 903.65216
$ws.Range("K77").Value = 4518.2608
$ws.Range("M77").Value = -150.2608

$ws.Range("H102").Value = 120418.336
$ws.Range("I102").Value = 120418.336
$ws.Range("K102").Value = 120418.336
$ws.Range("M102").Value = -118796.336

$ws.Range("H110").Value = 1999.4
$ws.Range("J110").Value = 999.3333
$ws.Range("L110").Value = 999.3333
$ws.Range("N110").Value = -5089.3333

$ws.Range("H122").Value = 2662.625
$ws.Range("I122").Value = 2627.3103
$ws.Range("J122").Value = 3004
$ws.Range("K122").Value = 7881.9309
$ws.Range("L122").Value = 9012
$ws.Range("M122").Value = -5431.9309
$ws.Range("N122").Value = -13912

$ws.Range("H132").Value = 1187930.8
$ws.Range("I132").Value = 1250243
$ws.Range("K132").Value = 3750729
$ws.Range("M132").Value = -3748199

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 564.75
$ws.Range("I11").Value = 91
$ws.Range("K11").Value = 91
$ws.Range("M11").Value = 49

$ws.Range("H94").Value = 1350.5
$ws.Range("J94").Value = 2998.5
$ws.Range("L94").Value = 2998.5
$ws.Range("N94").Value = -3900.5

$ws.Range("H99").Value = 10106.25
$ws.Range("I99").Value = 4708.6665
$ws.Range("J99").Value = 17046
$ws.Range("K99").Value = 4708.6665
$ws.Range("L99").Value = 17046
$ws.Range("M99").Value = -3210.6665
$ws.Range("N99").Value = -20042

$ws.Range("H107").Value = 8133139.5
$ws.Range("I107").Value = 3289.6177
$ws.Range("J107").Value = 47620984
$ws.Range("K107").Value = 3289.6177
$ws.Range("L107").Value = 47620984
$ws.Range("M107").Value = -1369.6177
$ws.Range("N107").Value = -47624824

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2016.8889
$ws.Range("I99").Value = 2100
$ws.Range("J99").Value = 2006.5
$ws.Range("K99").Value = 2100
$ws.Range("L99").Value = 2006.5
$ws.Range("M99").Value = -602
$ws.Range("N99").Value = -5002.5

$ws.Range("H122").Value = 3561.5
$ws.Range("I122").Value = 1991.3334
$ws.Range("J122").Value = 5131.6665
$ws.Range("K122").Value = 5974.0002
$ws.Range("L122").Value = 15394.9995
$ws.Range("M122").Value = -3524.0002
$ws.Range("N122").Value = -20294.9995

$ws.Range("H126").Value = 2016.8889
$ws.Range("I126").Value = 2100
$ws.Range("J126").Value = 2006.5
$ws.Range("K126").Value = 6300
$ws.Range("L126").Value = 6019.5
$ws.Range("M126").Value = -3830
$ws.Range("N126").Value = -10959.5

$ws.Range("H134").Value = 8349791.5
$ws.Range("I134").Value = 20517.572
$ws.Range("J134").Value = 37502250
$ws.Range("K134").Value = 61552.716
$ws.Range("L134").Value = 112506750
$ws.Range("M134").Value = -59017.716
$ws.Range("N134").Value = -112511820

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("M5").ClearContents()
$ws.Range("N5").ClearContents()

$ws.Range("H131").Value = 12815.571
$ws.Range("I131").Value = 1901.6666
$ws.Range("J131").Value = 14634.556
$ws.Range("K131").Value = 5704.9998
$ws.Range("L131").Value = 43903.66800000001
$ws.Range("M131").Value = -664.9997999999996
$ws.Range("N131").Value = -53983.66800000001

$ws.Range("H135").Value = 0
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("M135").ClearContents()
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7022.222
$ws.Range("I70").Value = 6900.125
$ws.Range("K70").Value = 6900.125
$ws.Range("M70").Value = -6630.125

$ws.Range("H73").Value = 7022.222
$ws.Range("I73").Value = 6900.125
$ws.Range("K73").Value = 6900.125
$ws.Range("M73").Value = -5964.125

$ws.Range("H97").Value = 3777.4443
$ws.Range("I97").Value = 3606.2666
$ws.Range("J97").Value = 4633.3335
$ws.Range("K97").Value = 3606.2666
$ws.Range("L97").Value = 4633.3335
$ws.Range("M97").Value = -3110.2666
$ws.Range("N97").Value = -5625.3335

$ws.Range("H102").Value = 3215.6135
$ws.Range("I102").Value = 2666.3845
$ws.Range("K102").Value = 2666.3845
$ws.Range("M102").Value = -1044.3845

$ws.Range("H132").Value = 126514376
$ws.Range("I132").Value = 168682270
$ws.Range("J132").Value = 10700
$ws.Range("K132").Value = 506046810
$ws.Range("L132").Value = 32100
$ws.Range("M132").Value = -506044280
$ws.Range("N132").Value = -37160

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 925.7059
$ws.Range("I16").Value = 917
$ws.Range("K16").Value = 917
$ws.Range("M16").Value = -747

$ws.Range("H68").Value = 4415
$ws.Range("I68").Value = 4766.222
$ws.Range("K68").Value = 4766.222
$ws.Range("M68").Value = -4017.222

$ws.Range("H71").Value = 4415
$ws.Range("I71").Value = 4766.222
$ws.Range("K71").Value = 23831.11
$ws.Range("M71").Value = -20087.11

$ws.Range("H93").Value = 913.5714
$ws.Range("I93").Value = 739
$ws.Range("K93").Value = 739
$ws.Range("M93").Value = 509

$ws.Range("H101").Value = 65750
$ws.Range("J101").Value = 65750
$ws.Range("L101").Value = 65750
$ws.Range("N101").Value = -72240

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 748.7778
$ws.Range("I126").Value = 748.7778
$ws.Range("K126").Value = 2246.3334
$ws.Range("M126").Value = 223.6666
